$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the old last row (row 9, "SVM") - the new table only has
#    7 data rows (LR, LDA, KNN, DTREE, RTREE, XTREE, SVM).
# ------------------------------------------------------------------
$ws.Rows("9").Delete()

# ------------------------------------------------------------------
# 2. Extend the header row with the new "mean"/"std" columns
#    (H1:L1), copying the existing header style (bold, centered,
#    bordered) from the existing G1 header cell so no new style
#    entries are introduced.
# ------------------------------------------------------------------
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("G1").Copy($ws.Range("I1"))
$ws.Range("G1").Copy($ws.Range("J1"))
$ws.Range("G1").Copy($ws.Range("K1"))
$ws.Range("G1").Copy($ws.Range("L1"))

# Rewrite every header label (C1:L1) - C1/D1 reuse the "One Year Alt"
# root name split into "mean"/"std" and so on for every horizon.
# New shared strings are appended in the order they're first written,
# so this order matches the expected sharedStrings.xml layout.
$ws.Range("C1").Value = "One Year Alt mean"
$ws.Range("D1").Value = "One Year Alt std"
$ws.Range("E1").Value = "Two Year Alt mean"
$ws.Range("F1").Value = "Two Year Alt std"
$ws.Range("G1").Value = "Three Year Alt mean"
$ws.Range("H1").Value = "Three Year Alt std"
$ws.Range("I1").Value = "Five Year Alt mean"
$ws.Range("J1").Value = "Five Year Alt std"
$ws.Range("K1").Value = "Ten Year Alt mean"
$ws.Range("L1").Value = "Ten Year Alt std"

# ------------------------------------------------------------------
# 3. Extend data rows 2:8 with the new H:L ("std") columns, copying
#    the (unstyled) number-cell formatting from an existing data
#    cell on the same row so no new styles are introduced.
# ------------------------------------------------------------------
foreach ($r in 2..8) {
    $ws.Range("G$r").Copy($ws.Range("H$r"))
    $ws.Range("G$r").Copy($ws.Range("I$r"))
    $ws.Range("G$r").Copy($ws.Range("J$r"))
    $ws.Range("G$r").Copy($ws.Range("K$r"))
    $ws.Range("G$r").Copy($ws.Range("L$r"))
}

# ------------------------------------------------------------------
# 4. Row 2 - LR
# ------------------------------------------------------------------
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8280016295994373
$ws.Range("D2").Value = 0.01314115558265945
$ws.Range("E2").Value = 0.8045483589888736
$ws.Range("F2").Value = 0.01021098054531338
$ws.Range("G2").Value = 0.7841161825352569
$ws.Range("H2").Value = 0.02037080918991764
$ws.Range("I2").Value = 0.7581045343232958
$ws.Range("J2").Value = 0.03686210095178553
$ws.Range("K2").Value = 0.7332453788165456
$ws.Range("L2").Value = 0.02100574089317048

# ------------------------------------------------------------------
# 5. Row 3 - LDA
# ------------------------------------------------------------------
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.828083341928125
$ws.Range("D3").Value = 0.01738954501765564
$ws.Range("E3").Value = 0.803975299753049
$ws.Range("F3").Value = 0.01878365712207545
$ws.Range("G3").Value = 0.7844136239200832
$ws.Range("H3").Value = 0.01662549635590544
$ws.Range("I3").Value = 0.7589713514487531
$ws.Range("J3").Value = 0.03626634152315474
$ws.Range("K3").Value = 0.7377752104952041
$ws.Range("L3").Value = 0.02606557454878471

# ------------------------------------------------------------------
# 6. Row 4 - KNN
# ------------------------------------------------------------------
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.811240568513268
$ws.Range("D4").Value = 0.01786549846085677
$ws.Range("E4").Value = 0.8063405954170133
$ws.Range("F4").Value = 0.01820939297161448
$ws.Range("G4").Value = 0.7959019063281368
$ws.Range("H4").Value = 0.01886414787937804
$ws.Range("I4").Value = 0.7912365362816405
$ws.Range("J4").Value = 0.02267233053982253
$ws.Range("K4").Value = 0.7842248902929911
$ws.Range("L4").Value = 0.02779921101222466

# ------------------------------------------------------------------
# 7. Row 5 - DTREE (was CART)
# ------------------------------------------------------------------
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7432808350473631
$ws.Range("D5").Value = 0.02970604243598193
$ws.Range("E5").Value = 0.7384687980288768
$ws.Range("F5").Value = 0.01461229853266952
$ws.Range("G5").Value = 0.7331453521949756
$ws.Range("H5").Value = 0.01595459853282057
$ws.Range("I5").Value = 0.7288539426334661
$ws.Range("J5").Value = 0.02764157540942656
$ws.Range("K5").Value = 0.6894208662034085
$ws.Range("L5").Value = 0.03890525283413538

# ------------------------------------------------------------------
# 8. Row 6 - RTREE
# ------------------------------------------------------------------
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.7462392054109436
$ws.Range("D6").Value = 0.01207996503351616
$ws.Range("E6").Value = 0.72807629247923
$ws.Range("F6").Value = 0.02251736131922065
$ws.Range("G6").Value = 0.7081733605484182
$ws.Range("H6").Value = 0.01599183527429442
$ws.Range("I6").Value = 0.6880605615572802
$ws.Range("J6").Value = 0.02223130337846224
$ws.Range("K6").Value = 0.6733541204691896
$ws.Range("L6").Value = 0.0187066950232615

# ------------------------------------------------------------------
# 9. Row 7 - XTREE
# ------------------------------------------------------------------
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8202843327676501
$ws.Range("D7").Value = 0.01368170429773109
$ws.Range("E7").Value = 0.8053424953033834
$ws.Range("F7").Value = 0.01480939687693497
$ws.Range("G7").Value = 0.7894694519649125
$ws.Range("H7").Value = 0.02412119819524255
$ws.Range("I7").Value = 0.7681696313779012
$ws.Range("J7").Value = 0.03052010661789345
$ws.Range("K7").Value = 0.7567199282853556
$ws.Range("L7").Value = 0.01994249456970832

# ------------------------------------------------------------------
# 10. Row 8 - SVM (was NB)
# ------------------------------------------------------------------
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8287908072786436
$ws.Range("D8").Value = 0.01381495469101302
$ws.Range("E8").Value = 0.8085283823907401
$ws.Range("F8").Value = 0.01664299585044512
$ws.Range("G8").Value = 0.8003995933194117
$ws.Range("H8").Value = 0.01762246748112882
$ws.Range("I8").Value = 0.7961138604671039
$ws.Range("J8").Value = 0.02876764975063725
$ws.Range("K8").Value = 0.7707720721685174
$ws.Range("L8").Value = 0.02301496505064199
